$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("A2").Value = "02e646fe-d018-41b1-81e0-88e810759ffe.md"
$ov.Range("B2").Value = "e2e\02e646fe-d018-41b1-81e0-88e810759ffe.md"
$ov.Hyperlinks.Item(1).TextToDisplay = "e2e\02e646fe-d018-41b1-81e0-88e810759ffe.md"
$ov.Range("G2").Value = "2016-08-23 17:04:50"

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("A2").Value = "02e646fe-d018-41b1-81e0-88e810759ffe.md"
$zh.Range("G2").Value = "02e646fe-d018-41b1-81e0-88e810759ffe.ed5a636d5288526de773b3633ea6b651012cee2d.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-23 17:04:46"
$zh.Range("I2").Value = ""
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = "0001-01-01 00:00:00"
$zh.Hyperlinks.Item(1).TextToDisplay = "02e646fe-d018-41b1-81e0-88e810759ffe.md"
$zh.Hyperlinks.Item(2).Delete()

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("A2").Value = "02e646fe-d018-41b1-81e0-88e810759ffe.md"
$de.Range("G2").Value = "02e646fe-d018-41b1-81e0-88e810759ffe.ed5a636d5288526de773b3633ea6b651012cee2d.de-de.xlf"
$de.Range("H2").Value = "2016-08-23 17:04:50"
$de.Range("I2").Value = ""
$de.Range("J2").Value = ""
$de.Range("K2").Value = "0001-01-01 00:00:00"
$de.Hyperlinks.Item(1).TextToDisplay = "02e646fe-d018-41b1-81e0-88e810759ffe.md"
$de.Hyperlinks.Item(2).Delete()
